# Registers.xlsx update — RiskRegister sheet, row 12 ("2 Challenging tasks" risk)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RiskRegister")

# Fill in the previously-empty risk row (row 12) with the new risk entry.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "2 Challenging tasks"
$ws.Range("C12").Value = "Too much work, one of the two could not work out"
$ws.Range("D12").Value = "Final deliverable not reached"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = "Can not use previous years' work"
$ws.Range("G12").Value = 2
$ws.Range("L12").Value = "Work in parallell, make subteams"

# Leave H12/I12 blank and K12 at its existing 0 (unchanged); J12 is the
# pre-existing shared formula (=E12*G12) and recalculates automatically.

# Move the active selection to B12, matching the author's final cursor spot.
$ws.Range("B12").Select()
